# Applies the cover-letter.docx revision described by the target diff.
#
# The edit mixes three kinds of changes inside existing paragraphs:
#   1) adjacent runs with identical formatting get merged into one run
#      (no visible text change) in two places in the opening paragraph;
#   2) a proofErr spellStart/spellEnd wrapper around the word "the " is
#      removed and its run re-merged with its neighbours (no visible text
#      change) in the results paragraph;
#   3) a large new block of (partly italicised) citation text replaces
#      part of the "journal selection" paragraph;
#   4) a <w:lastRenderedPageBreak/> marker is added to the run that starts
#      the closing paragraph.
#
# None of these can be expressed with plain Find/Replace because they
# involve precise run/formatting boundaries (italics, proofErr wrappers).
# Instead, for each affected paragraph we take the Range spanning its
# text (excluding the trailing paragraph mark, i.e. Start .. End-1) and
# call Range.InsertXML with a full WordProcessingML package whose body
# holds the desired replacement runs; InsertXML splices that content in
# place of the selected range while leaving the paragraph's own mark (and
# therefore its paragraph-level identity/formatting/section link) intact.

$d = $word.ActiveDocument

function Replace-ParagraphRuns($paragraphIndex, $innerXml) {
    $full = $d.Paragraphs($paragraphIndex).Range
    $target = $d.Range($full.Start, $full.End - 1)
    $target.InsertXML($innerXml)
}

# --- Paragraph 6: "Thank you for considering our article titled, ..." ---
# Two pairs of adjacent runs get merged into single runs; the rendered text
# is unchanged.
$xmlPara6 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:r><w:t xml:space="preserve">Thank you for considering our article titled, “Abiotic and biotic factors jointly influence the transmission of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ranavirus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in larval amphibian communities” for publication in </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Functional Ecology</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphRuns 6 $xmlPara6

# --- Paragraph 8: "We study the joint influence of biotic ..." ---
# The spellStart/"the "/spellEnd split around "the host species" collapses
# into the surrounding run; rendered text is unchanged.
$xmlPara8 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:r><w:t xml:space="preserve">We study the joint influence of biotic and abiotic factors by focusing on a limited number of factors and the mechanisms underlying their effects on multimodal transmission of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ranavirus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in larval amphibian communities. By focusing on host community composition, host abundance, and the environmental persistence of the pathogen, we developed a simple mathematical model that can examine each of these factors independently and together. Through our model, we learn that each of these factors can contribute significantly to transmission and that the joint influence of these factors can have synergistic effects. We connect this model to empirical data </w:t></w:r><w:r w:rsidR="002864D8"><w:t xml:space="preserve">of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="002864D8"><w:t>ranavirus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="002864D8"><w:t xml:space="preserve"> infections in linked communities of amphibians and find that each of these factors can overlap and appear to contribute to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="002864D8"><w:t>ranavirus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="002864D8"><w:t xml:space="preserve"> prevalence. When host abundance is high, communities are typically dominated by high competence hosts and this typically occurs in cooler months, when lower water temperatures may enhance viral persistence times in the water. Furthermore, we examine the phylogenetic relationship of the host species and find that higher values of competence are moderately spread throughout the phylogeny, potentially enabling coexistence of high competence species. The consideration of both biotic and abiotic factors and the mechanistic rationale that we provide can contribute significantly to our understanding of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="002864D8"><w:t>ranavirus</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="002864D8"><w:t xml:space="preserve"> epidemics in natural systems and can also further our understanding of other </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="002864D8"><w:t>multimodel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="002864D8"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="002864D8"><w:t>multihost</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="002864D8"><w:t xml:space="preserve"> pathogens. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphRuns 8 $xmlPara8

# --- Paragraph 10: "We have chosen to submit out article to ..." ---
# The tail of the paragraph is rewritten to cite Rohr et al. 2020, Shaw &
# Civitello 2021, Dobson 2004 and Johnson et al. 2013, replacing the old
# "diversity-disease literature has been dominated by debate ..." passage.
$xmlPara10 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:r><w:t xml:space="preserve">We have chosen to submit out article to </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Functional Ecology</w:t></w:r><w:r><w:t xml:space="preserve"> because of the journal’s strong record of published research regarding biodiversity-disease relationships and the mechanisms underlying these relationships. </w:t></w:r><w:r><w:t>The study of the community ecology of infectious diseases has already benefitted from analyses that focus on developing a mechanistic understanding of biodiversity-disease relationships (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>e.g.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Rohr et al. 2020, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Nature Ecol</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve"> &amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Evol</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> 4, 24-33; Shaw &amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Civitello</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 2021, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Functional Ecology</w:t></w:r><w:r><w:t xml:space="preserve"> 35, 2376-2386). </w:t></w:r><w:r><w:t xml:space="preserve">We intend to contribute to this catalog by expanding theory that focuses on </w:t></w:r><w:r><w:t xml:space="preserve">the specific compositions of communities and how these can drive transmission in conjunction with other important abiotic factors. </w:t></w:r><w:r><w:t xml:space="preserve">We follow recommendations from </w:t></w:r><w:r><w:t xml:space="preserve">Shaw &amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Civitello</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> (2021) </w:t></w:r><w:r><w:t xml:space="preserve">and focus on ecological interactions rather than species richness and use the parasite’s basic reproductive ratio at the community scale to assess our models (Dobson 2004, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Am. Nat.</w:t></w:r><w:r><w:t xml:space="preserve"> 164, S64-S78). We also attempt to address research frontiers outlined in Rohr et al. (2020) by studying the context dependencies that may explain our ability to detect biodiversity-disease patterns and use other metrics such as community competence (</w:t></w:r><w:r><w:t xml:space="preserve">Johnson et al. 2013, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>Nature</w:t></w:r><w:r><w:t xml:space="preserve"> 494, 230-233</w:t></w:r><w:r><w:t>) which can further help</w:t></w:r><w:r><w:t xml:space="preserve"> us</w:t></w:r><w:r><w:t xml:space="preserve"> to understand the processes driving the patterns we observe</w:t></w:r><w:r><w:t>d</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> Our study, along with many others, are trending towards developing theory for understanding how abiotic factors influence community composition and how these formulations of host communities relate to disease transmission. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphRuns 10 $xmlPara10

# --- Paragraph 12: "This paper has not been published anywhere else ..." ---
# A <w:lastRenderedPageBreak/> marker is added at the start of the run.
$xmlPara12 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">This paper has not been published anywhere else and is not under consideration at any other journals. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
Replace-ParagraphRuns 12 $xmlPara12
